$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$ws.Range("B3").Value = "6.0.0"

# Update Date value (row 8)
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Update Publisher value (row 9) - was empty
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it entirely,
# shifting Description/Purpose/Copyright/Immutable rows up by one.
$ws.Rows.Item(11).Delete()
